$wb = $excel.ActiveWorkbook

# Rename InvalidLogin -> Dashboard, add Program sheet
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "Dashboard"

$ws3 = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Program"

# ---- Dashboard sheet data ----
$ws2.Range("A1").Value = "KeyOption"
$ws2.Range("B1").Value = "UserNameData"
$ws2.Range("C1").Value = "PasswordData"
$ws2.Range("D1").Value = "expectedValue"
$ws2.Range("E1").Value = "nth"

$ws2.Range("A2").Value = "title"
$ws2.Range("B2").Value = "playwrightuser@gmail.com"
$ws2.Range("C2").Value = "Playwright@1234"
$ws2.Range("D2").Value = " LMS - Learning Management System "

$ws2.Range("A3").Value = "navBarText"
$ws2.Range("B3").Value = "playwrightuser@gmail.com"
$ws2.Range("C3").Value = "Playwright@1234"
$ws2.Range("D3").Value = " LMS - Learning Management System ProgramBatchClassLogout"

$ws2.Range("A4").Value = "program"
$ws2.Range("B4").Value = "playwrightuser@gmail.com"
$ws2.Range("C4").Value = "Playwright@1234"
$ws2.Range("D4").Value = "Program"
$ws2.Range("E4").Value = 0

$ws2.Range("A5").Value = "batch"
$ws2.Range("B5").Value = "playwrightuser@gmail.com"
$ws2.Range("C5").Value = "Playwright@1234"
$ws2.Range("D5").Value = "Batch"
$ws2.Range("E5").Value = 1

$ws2.Range("A6").Value = "class"
$ws2.Range("B6").Value = "playwrightuser@gmail.com"
$ws2.Range("C6").Value = "Playwright@1234"
$ws2.Range("D6").Value = "Class"
$ws2.Range("E6").Value = 2

$ws2.Range("A7").Value = "logout"
$ws2.Range("B7").Value = "playwrightuser@gmail.com"
$ws2.Range("C7").Value = "Playwright@1234"
$ws2.Range("D7").Value = "Logout"
$ws2.Range("E7").Value = 3

# ---- Program sheet data ----
$ws3.Range("A1").Value = "KeyOption"
$ws3.Range("B1").Value = "UserNameData"
$ws3.Range("C1").Value = "PasswordData"
$ws3.Range("D1").Value = "expected"

$ws3.Range("A2").Value = "programModule"
$ws3.Range("B2").Value = "playwrightuser@gmail.com"
$ws3.Range("C2").Value = "Playwright@1234"
$ws3.Range("D2").Value = "Manage Program"

$ws3.Range("B3").Value = "playwrightuser@gmail.com"
$ws3.Range("C3").Value = "Playwright@1234"

$ws3.Range("B4").Value = "playwrightuser@gmail.com"
$ws3.Range("C4").Value = "Playwright@1234"

$ws3.Range("B5").Value = "playwrightuser@gmail.com"
$ws3.Range("C5").Value = "Playwright@1234"

$ws3.Range("B6").Value = "playwrightuser@gmail.com"
$ws3.Range("C6").Value = "Playwright@1234"

$ws3.Range("B7").Value = "playwrightuser@gmail.com"
$ws3.Range("C7").Value = "Playwright@1234"

# ---- Sheet views / selections ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Select()
$excel.ActiveWindow.SelectedSheets.Item(1).Range("G5").Select()

$ws2.Select()
$ws2.Range("A7").Select()

$ws3.Select()
$ws3.Range("G12").Select()

$ws2.Select()

Write-Host "done"
